$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.207534670829773
$ws.Range("B1").Value = 2.622440814971924
$ws.Range("D1").Value = 2.167834758758545
$ws.Range("E1").Value = 1.165008068084717
